# edit.ps1 - Apply Phantom_Profits cell-value corrections across all 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Values derived from the supplied OOXML diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (54 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 312.57144
$ws.Range("I2").Value = 306.33334
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 306.33334
$ws.Range("L2").Value = 350
$ws.Range("M2").Value = -193.33334
$ws.Range("N2").Value = -576
# row 17
$ws.Range("H17").Value = 2997.3333
$ws.Range("J17").Value = 2997.3333
$ws.Range("L17").Value = 8991.999899999999
$ws.Range("N17").Value = -9327.999899999999
# row 19
$ws.Range("H19").Value = 909.6
$ws.Range("I19").Value = 959.36365
$ws.Range("J19").Value = 772.75
$ws.Range("K19").Value = 959.36365
$ws.Range("L19").Value = 772.75
$ws.Range("M19").Value = -784.36365
$ws.Range("N19").Value = -1122.75
# row 29
$ws.Range("H29").Value = 24875
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 33000
$ws.Range("K29").Value = 1500
$ws.Range("L29").Value = 99000
$ws.Range("M29").Value = -1219
$ws.Range("N29").Value = -99562
# row 38
$ws.Range("H38").Value = 391
$ws.Range("I38").Value = 301.1111
$ws.Range("J38").Value = 1200
$ws.Range("K38").Value = 903.3333
$ws.Range("L38").Value = 3600
$ws.Range("M38").Value = -531.3333
$ws.Range("N38").Value = -4344
# row 47
$ws.Range("H47").Value = 74000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 74000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 74000
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -75944
# row 58
$ws.Range("H58").Value = 2888.75
$ws.Range("I58").Value = 527.5
$ws.Range("J58").Value = 5250
$ws.Range("K58").Value = 1582.5
$ws.Range("L58").Value = 15750
$ws.Range("M58").Value = -1432.5
$ws.Range("N58").Value = -16050
# row 129
$ws.Range("H129").Value = 2809.5557
$ws.Range("I129").Value = 1947.75
$ws.Range("K129").Value = 5843.25
$ws.Range("M129").Value = -843.25
# row 132
$ws.Range("H132").Value = 6332.0527
$ws.Range("I132").Value = 6947.3076
$ws.Range("K132").Value = 20841.9228
$ws.Range("M132").Value = -18311.9228

# ---- Sheet: ARM (18 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
# row 88
$ws.Range("H88").Value = 2495.9333
$ws.Range("I88").Value = 2314.7144
$ws.Range("J88").Value = 2654.5
$ws.Range("K88").Value = 2314.7144
$ws.Range("L88").Value = 2654.5
$ws.Range("M88").Value = -1908.7144
$ws.Range("N88").Value = -3466.5
# row 91
$ws.Range("H91").Value = 2495.9333
$ws.Range("I91").Value = 2314.7144
$ws.Range("J91").Value = 2654.5
$ws.Range("K91").Value = 2314.7144
$ws.Range("L91").Value = 2654.5
$ws.Range("M91").Value = -910.7143999999998
$ws.Range("N91").Value = -5462.5
# row 132
$ws.Range("H132").Value = 2729.5264
$ws.Range("I132").Value = 2729.5264
$ws.Range("K132").Value = 8188.5792
$ws.Range("M132").Value = -5658.5792

# ---- Sheet: BSM (23 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
# row 11
$ws.Range("H11").Value = 5456
$ws.Range("J11").Value = 8000
$ws.Range("L11").Value = 8000
$ws.Range("N11").Value = -8280
# row 99
$ws.Range("H99").Value = 2832
$ws.Range("I99").Value = 1200
$ws.Range("J99").Value = 3158.4
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 3158.4
$ws.Range("M99").Value = 298
$ws.Range("N99").Value = -6154.4
# row 100
$ws.Range("H100").Value = 23574.75
$ws.Range("J100").Value = 23574.75
$ws.Range("L100").Value = 23574.75
$ws.Range("N100").Value = -25738.75
# row 105
$ws.Range("H105").Value = 2918.6
$ws.Range("I105").Value = 2898.25
$ws.Range("K105").Value = 2898.25
$ws.Range("M105").Value = -1151.25
# row 107
$ws.Range("H107").Value = 2164.1538
$ws.Range("I107").Value = 2164.1538
$ws.Range("K107").Value = 2164.1538
$ws.Range("M107").Value = -244.1538

# ---- Sheet: CRP (19 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("I22").Value = 2024.375
$ws.Range("K22").Value = 2024.375
$ws.Range("M22").Value = -1674.375
# row 31
$ws.Range("H31").Value = 3054.5334
$ws.Range("I31").Value = 2219.375
$ws.Range("K31").Value = 2219.375
$ws.Range("M31").Value = -1924.375
# row 34
$ws.Range("H34").Value = 3054.5334
$ws.Range("I34").Value = 2219.375
$ws.Range("K34").Value = 2219.375
$ws.Range("M34").Value = -2017.375
# row 118
$ws.Range("H118").Value = 42999.6
$ws.Range("J118").Value = 42999.6
$ws.Range("L118").Value = 42999.6
$ws.Range("N118").Value = -46313.6
# row 132
$ws.Range("H132").Value = 3246.0625
$ws.Range("I132").Value = 2869.75
$ws.Range("K132").Value = 8609.25
$ws.Range("M132").Value = -6079.25

# ---- Sheet: CUL (15 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
# row 17
$ws.Range("H17").Value = 2298.889
$ws.Range("I17").Value = 170
$ws.Range("J17").Value = 2907.1428
$ws.Range("K17").Value = 510
$ws.Range("L17").Value = 8721.428400000001
$ws.Range("M17").Value = -341
$ws.Range("N17").Value = -9059.428400000001
# row 39
$ws.Range("H39").Value = 6623.3
$ws.Range("J39").Value = 6525.8887
$ws.Range("L39").Value = 19577.6661
$ws.Range("N39").Value = -20165.6661
# row 55
$ws.Range("H55").Value = 401
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

# ---- Sheet: GSM (4 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
# row 122
$ws.Range("H122").Value = 1919.25
$ws.Range("I122").Value = 1919.25
$ws.Range("K122").Value = 5757.75
$ws.Range("M122").Value = -3307.75

# ---- Sheet: LTW (32 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 1164.6666
$ws.Range("I22").Value = 1164.6666
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1164.6666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -869.6666
$ws.Range("N22").ClearContents()
# row 27
$ws.Range("H27").Value = 1164.6666
$ws.Range("I27").Value = 1164.6666
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1164.6666
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -1057.6666
$ws.Range("N27").ClearContents()
# row 46
$ws.Range("H46").Value = 1586
$ws.Range("I46").Value = 1149.909
$ws.Range("J46").Value = 3185
$ws.Range("K46").Value = 1149.909
$ws.Range("L46").Value = 3185
$ws.Range("M46").Value = -961.9090000000001
$ws.Range("N46").Value = -3561
# row 55
$ws.Range("H55").Value = 1079.6428
$ws.Range("I55").Value = 352.83334
$ws.Range("J55").Value = 1624.75
$ws.Range("K55").Value = 352.83334
$ws.Range("L55").Value = 1624.75
$ws.Range("M55").Value = -179.83334
$ws.Range("N55").Value = -1970.75
# row 98
$ws.Range("H98").Value = 59998.2
$ws.Range("J98").Value = 59998.2
$ws.Range("L98").Value = 59998.2
$ws.Range("N98").Value = -65988.2

# ---- Sheet: WVR (32 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
# row 81
$ws.Range("H81").Value = 2137.6365
$ws.Range("I81").Value = 2051.4
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 4102.8
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -3041.8
$ws.Range("N81").Value = -8122
# row 84
$ws.Range("H84").Value = 2137.6365
$ws.Range("I84").Value = 2051.4
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 20514
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -15210
$ws.Range("N84").Value = -40608
# row 126
$ws.Range("H126").Value = 4845.7
$ws.Range("I126").Value = 4884.1113
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 14652.3339
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -12182.3339
$ws.Range("N126").Value = -18440
# row 130
$ws.Range("H130").Value = 84950
$ws.Range("J130").Value = 84950
$ws.Range("L130").Value = 84950
$ws.Range("N130").Value = -94990
# row 132
$ws.Range("H132").Value = 3323.5
$ws.Range("I132").Value = 2738.7
$ws.Range("J132").Value = 6247.5
$ws.Range("K132").Value = 8216.099999999999
$ws.Range("L132").Value = 18742.5
$ws.Range("M132").Value = -5686.099999999999
$ws.Range("N132").Value = -23802.5

